$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.490.45"
$ws.Range("E2").Value = "  +4.15%  "

$ws.Range("D3").Value = "2.626.20"
$ws.Range("E3").Value = "  +4.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.52"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.03"
$ws.Range("E6").Value = "  +1.42%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  +1.67%  "

$ws.Range("E9").Value = "  +9.96%  "

$ws.Range("D10").Value = "2.625.11"
$ws.Range("E10").Value = "  +4.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("E12").Value = "  +2.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.04"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000189"
$ws.Range("E14").Value = "  +3.49%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.104.75"
$ws.Range("E15").Value = "  +4.08%  "

$ws.Range("D16").Value = "72.322.32"
$ws.Range("E16").Value = "  +3.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.66"
$ws.Range("E17").Value = "  +2.45%  "

$ws.Range("D18").Value = "2.623.37"
$ws.Range("E18").Value = "  +5.45%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.00"
$ws.Range("E19").Value = "  +4.74%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "380.77"
$ws.Range("E20").Value = "  +4.63%  "

$ws.Range("E21").Value = "  +5.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("E22").Value = "  +3.43%  "

$ws.Range("E23").Value = "  +18.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.19"
$ws.Range("E24").Value = "  +3.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.39"
$ws.Range("E26").Value = "  +3.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  +9.19%  "

$ws.Range("D28").Value = "2.759.93"
$ws.Range("E28").Value = "  +4.14%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "0.0₃0957"
$ws.Range("E30").Value = "  +6.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "520.40"
$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  +4.36%  "

$ws.Range("E33").Value = "  +6.77%  "

$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.18"
$ws.Range("E36").Value = "  +3.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.34"
$ws.Range("E37").Value = "  +3.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.08"
$ws.Range("E38").Value = "  +0.92%  "

$ws.Range("E39").Value = "  +6.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.110"
$ws.Range("E40").Value = "  -8.00%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.08"
$ws.Range("E42").Value = "  +5.75%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.59"
$ws.Range("E44").Value = "  +8.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.331"
$ws.Range("E45").Value = "  +2.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.52"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.04"
$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.71"
$ws.Range("E48").Value = "  +3.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.544"
$ws.Range("E49").Value = "  +5.17%  "

$ws.Range("E50").Value = "  +7.21%  "

$ws.Range("D51").Value = "0.0₆0263"
$ws.Range("E51").Value = "  +4.89%  "
